$d = $word.ActiveDocument

# Replacement 1: 'Processing of Ceramics I'
$d.Content.Find.Execute("Processing of Ceramics I", $true, $false, $false, $false, $false, $true, 1, $false, "Processing of Ceramics", 2) | Out-Null

# Replacement 2: 'Créditos-trabalho: 1'
$d.Content.Find.Execute("Créditos-trabalho: 1", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-trabalho: 0", 2) | Out-Null

# Replacement 3: 'Carga horária: 90 h'
$d.Content.Find.Execute("Carga horária: 90 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 60 h", 2) | Out-Null

# Replacement 4: 'Ativação: 01/01/2023'
$d.Content.Find.Execute("Ativação: 01/01/2023", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# Replacement 5: '1 – Introdução, 2 – Principais matérias-primas na
$d.Content.Find.Execute("1 – Introdução, 2 – Principais matérias-primas naturais e sintéticas, 3 - Preparo de massas cerâmicas, 4 – Conformação, 5 – Queima (sintereização), 6 - Variáveis críticas no controle do processamento, 7 – Usinagem de materiais cerâmicos. 8 – Projeto", $true, $false, $false, $false, $false, $true, 1, $false, "1 – Introdução, 2 – Principais matérias-primas naturais e sintéticas, 3 - Preparo de massas cerâmicas, 4 – Conformação, 5 – Queima (sinterização), 6 – Projeto", 2) | Out-Null

# Replacement 6: '1 – Introdução: definições, setores cerâmicos, ap
$d.Content.Find.Execute("1 – Introdução: definições, setores cerâmicos, aplicações e fontes de divulgação da área de cerâmicas. 2 - Principais matérias-primas naturais e sintéticas: ocorrencias naturais e beneficiamento e sínteses de cerâmicas avançadas. 3 - Preparo de massas cerâmicas: formulação de composições cerâmicas com e sem utilização de diagramas de fases, reologia das barbotinas e pastas. 4 – Conformação:  equipamentos utilizados na conformação de cerâmicas tradicionais e técnicas, defeitos e problemas na conformação, métodos de conformação (colagem de barbotina, prensagem, extrusão, injeção). 5 – Queima (sintereização): curva de queima, eventos pré-sinterização, sinterização, mecanismos de sinterização, equipamentos, sinterização rápida, microestrutura (controle microestrutural, relação microestrutura x propriedades), 6 - Variáveis críticas no controle do processamento: avaliadas em cada etapa do processamento. 7 - Usinagem de materiais cerâmicos: usinagem a verde e após sinterização, defeitos superficiais introduzidos, acabamento.8 – Projeto: Desenvolvimento de produtos cerâmicos levando em conta aspectos de inovação, sustentabilidade, social e ecônômico. Este tópico deverá ser desenvolvido em grupo.", $true, $false, $false, $false, $false, $true, 1, $false, "1 – Introdução: definições, setores cerâmicos, aplicações e fontes de divulgação da área de cerâmicas. 2 - Principais matérias-primas naturais e sintéticas: ocorrências naturais e beneficiamento e sínteses de cerâmicas avançadas. 3 - Preparo de massas cerâmicas: moagem, análise granulométrica, formulação de composições cerâmicas com e sem utilização de diagramas de fases, reologia das barbotinas e pastas. 4 – Conformação:  equipamentos utilizados na conformação de cerâmicas tradicionais e técnicas, defeitos e problemas na conformação, métodos de conformação (colagem de barbotina, prensagem, extrusão convencional e 3D (manufatura aditiva), injeção). 5 – Queima (sinterização): curva de queima, eventos pré-sinterização, sinterização, mecanismos de sinterização, equipamentos, sinterização rápida, microestrutura (controle microestrutural, relação microestrutura x propriedades), 6 – Projeto: Desenvolvimento de produtos cerâmicos levando em conta aspectos de inovação, sustentabilidade, social e econômico. Este tópico deverá ser desenvolvido em grupo.", 2) | Out-Null
